# Natmi following Dr Hou advice
# Recomputes the Ccl2->Cxcr3 LR-pair table: each "Sending cluster" now gets a
# row per "Target cluster" (M1 and M2), replacing the old single-target-only
# rows, and all of the derived statistics are refreshed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..T for data rows 2..13 (header row 1 is untouched)
$data = @(
    @("ECs",    "Ccl2", "Cxcr3", "M1", 2, 1, 7.115595,           14.23119,           0.009306976103932288, 0.006258695489408662, 2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 5.099739820905,     30.59843892543,      0.002118461149705837, 0.001424609142012242),
    @("ECs",    "Ccl2", "Cxcr3", "M2", 2, 1, 7.115595,           14.23119,           0.009306976103932288, 0.006258695489408662, 3, 1,                  2.431954666666666,  7.295864, 0.7723792211295388, 0.7723792211295388, 17.30480446636,     103.82882679816,     0.00718851495422645,  0.00483408634739642),
    @("FAPs",   "Ccl2", "Cxcr3", "M1", 3, 1, 392.0055033333333,  1176.01651,         0.5127309595519816,   0.5171970317736687,   2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 280.9499522334967,  2528.54957010147,    0.116708220364221,    0.1177247912018131),
    @("FAPs",   "Ccl2", "Cxcr3", "M2", 3, 1, 392.0055033333333,  1176.01651,         0.5127309595519816,   0.5171970317736687,   3, 1,                  2.431954666666666,  7.295864, 0.7723792211295388, 0.7723792211295388, 953.3396131905155,  8580.056518714639,   0.3960227391877606,   0.3994722405718556),
    @("M1",     "Ccl2", "Cxcr3", "M1", 3, 1, 225.9980216666667,  677.9940650000001, 0.2955983564533449,   0.2981731251189355,   2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 161.9725561304784,  1457.753005174305,   0.06728432812873857,  0.06787039897781157),
    @("M1",     "Ccl2", "Cxcr3", "M2", 3, 1, 225.9980216666667,  677.9940650000001, 0.2955983564533449,   0.2981731251189355,   3, 1,                  2.431954666666666,  7.295864, 0.7723792211295388, 0.7723792211295388, 549.6169434496844,  4946.55249104716,    0.2283140283246063,   0.2303027261411239),
    @("M2",     "Ccl2", "Cxcr3", "M1", 3, 1, 119.4232813333333,  358.269844,         0.1562019234684543,   0.1575624987533681,   2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 85.59054630831866,  770.314916774868,    0.03555480348095373,  0.03586449868701772),
    @("M2",     "Ccl2", "Cxcr3", "M2", 3, 1, 119.4232813333333,  358.269844,         0.1562019234684543,   0.1575624987533681,   3, 1,                  2.431954666666666,  7.295864, 0.7723792211295388, 0.7723792211295388, 290.4320063472462,  2613.888057125216,   0.1206471199875006,   0.1216980000663504),
    @("Neutro", "Ccl2", "Cxcr3", "M1", 3, 1, 7.311581666666666,  21.934745,          0.009563320545030469, 0.009646620528067518, 2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 5.240203268918334,  47.161829420265,     0.002176810471047719, 0.002195771278066508),
    @("Neutro", "Ccl2", "Cxcr3", "M2", 3, 1, 7.311581666666666,  21.934745,          0.009563320545030469, 0.009646620528067518, 3, 1,                  2.431954666666666,  7.295864, 0.7723792211295388, 0.7723792211295388, 17.78143515496444,  160.03291639468,     0.00738651007398275,  0.00745084925000101),
    @("sCs",    "Ccl2", "Cxcr3", "M1", 2, 1, 12.69026,           25.38052,           0.01659846387725661,  0.01116202833655136,  2, 0.6666666666666666, 0.7166990000000001, 2.150097, 0.2276207788704612, 0.2276207788704611, 9.09509665174,      54.57057991044,      0.003778155275794364, 0.002540709583739978),
    @("sCs",    "Ccl2", "Cxcr3", "M2", 2, 1, 12.69026,           25.38052,           0.01659846387725661,  0.01116202833655136,  3, 1,                  2.431954666666666,  7.295864, 0.7723792211295388, 0.7723792211295388, 30.86213702821333,  185.17282216928,     0.01282030860146224,  0.008621318752811379)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = 2 + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = 1 + $j
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
